$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the data ranges of D and E columns (rows 2-51) so that
# numeric-looking strings (e.g. "4.08", "0.506") are stored as text, matching
# the original inlineStr/text cell contents, rather than being auto-converted
# to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.298.91'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '1.589.97'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").Value = '211.09'
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").Value = '0.506'
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("E8").Value = '  +0.35%  '
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("D12").Value = '1.816.00'
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.600.03'
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '4.08'
$ws.Range("E14").Value = '  +1.30%  '
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("D16").Value = '64.61'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").Value = '26.304.82'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '0.0₃0731'
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("D19").Value = '7.51'
$ws.Range("E19").Value = '  +3.86%  '
$ws.Range("D20").Value = '212.22'
$ws.Range("E20").Value = '  +2.38%  '
$ws.Range("E21").Value = '  -0.36%  '
$ws.Range("D22").Value = '4.29'
$ws.Range("E22").Value = '  +0.38%  '
$ws.Range("D23").Value = '9.02'
$ws.Range("E23").Value = '  +2.04%  '
$ws.Range("E24").Value = '  -2.88%  '
$ws.Range("D25").Value = '143.90'
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("E27").Value = '  +0.96%  '
$ws.Range("E28").Value = '  -0.83%  '
$ws.Range("E29").Value = '  -0.45%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("D33").Value = '2.99'
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("D34").Value = '1.325.25'
$ws.Range("E34").Value = '  +3.21%  '
$ws.Range("E35").Value = '  -1.60%  '
$ws.Range("D36").Value = '0.603'
$ws.Range("E36").Value = '  -0.74%  '
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("D39").Value = '0.816'
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("E40").Value = '  -0.37%  '
$ws.Range("E41").Value = '  +4.74%  '
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -23.82%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").Value = '0.765'
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.727.61'
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '61.91'
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("D47").Value = '88.04'
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("E48").Value = '  -5.08%  '
$ws.Range("E49").Value = '  -0.97%  '
$ws.Range("D50").Value = '0.0979'
$ws.Range("E50").Value = '  -4.52%  '
$ws.Range("E51").Value = '  -0.30%  '

# Restore the default (unstyled) cell style on the data ranges so the
# text-format coercion above does not leave a stray style index on cells
# that originally had no explicit style.
$ws.Range("D2:E51").Style = "Normal"
